$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (before the former row 89),
# shifting all existing rows 89-108 down to 91-110.
$ws.Rows.Item(89).Insert()
$ws.Rows.Item(89).Insert()

# New row 89: latest "Magnum" quote
$ws.Range("A89").Value = 2
$ws.Range("B89").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C89").Value = "Coquimbo"
$ws.Range("D89").Value = 44504
$ws.Range("D89").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E89").Value = 4
$ws.Range("F89").Value = 100112031
$ws.Range("G89").Value = "Poroto verde"
$ws.Range("H89").Value = "Magnum"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 34000
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = 34500
$ws.Range("N89").Value = "$/malla 25 kilos"
$ws.Range("O89").Value = "Provincia de Limarí"
$ws.Range("P89").Value = 1380
$ws.Range("Q89").Value = 25
$ws.Range("R89").Value = "Hortaliza"

# New row 90: latest "Sin especificar" quote
$ws.Range("A90").Value = 2
$ws.Range("B90").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C90").Value = "Coquimbo"
$ws.Range("D90").Value = 44504
$ws.Range("D90").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E90").Value = 4
$ws.Range("F90").Value = 100112031
$ws.Range("G90").Value = "Poroto verde"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 360
$ws.Range("K90").Value = 40000
$ws.Range("L90").Value = 42000
$ws.Range("M90").Value = 41000
$ws.Range("N90").Value = "$/malla 25 kilos"
$ws.Range("O90").Value = "Provincia de Limarí"
$ws.Range("P90").Value = 1640
$ws.Range("Q90").Value = 25
$ws.Range("R90").Value = "Hortaliza"
